$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4407
$ws.Range("I62").Value = 3954.7144
$ws.Range("J62").Value = 4859.2856
$ws.Range("K62").Value = 3954.7144
$ws.Range("L62").Value = 4859.2856
$ws.Range("M62").Value = -3330.7144
$ws.Range("N62").Value = -6107.2856
$ws.Range("H64").Value = 3880.5066
$ws.Range("I64").Value = 3799.9714
$ws.Range("J64").Value = 3947.6191
$ws.Range("K64").Value = 3799.9714
$ws.Range("L64").Value = 3947.6191
$ws.Range("M64").Value = -3551.9714
$ws.Range("N64").Value = -4443.6191
$ws.Range("H65").Value = 4407
$ws.Range("I65").Value = 3954.7144
$ws.Range("J65").Value = 4859.2856
$ws.Range("K65").Value = 19773.572
$ws.Range("L65").Value = 24296.428
$ws.Range("M65").Value = -16653.572
$ws.Range("N65").Value = -30536.428
$ws.Range("H67").Value = 3880.5066
$ws.Range("I67").Value = 3799.9714
$ws.Range("J67").Value = 3947.6191
$ws.Range("K67").Value = 3799.9714
$ws.Range("L67").Value = 3947.6191
$ws.Range("M67").Value = -2941.9714
$ws.Range("N67").Value = -5663.6191
$ws.Range("H111").Value = 4519.1
$ws.Range("I111").Value = 5822.7144
$ws.Range("J111").Value = 1477.3334
$ws.Range("K111").Value = 17468.1432
$ws.Range("L111").Value = 4432.0002
$ws.Range("M111").Value = -14401.1432
$ws.Range("N111").Value = -10566.0002
$ws.Range("H113").Value = 1624.75
$ws.Range("I113").Value = 1080
$ws.Range("J113").Value = 1806.3334
$ws.Range("K113").Value = 1080
$ws.Range("L113").Value = 1806.3334
$ws.Range("M113").Value = 2174
$ws.Range("N113").Value = -8314.3334
$ws.Range("H116").Value = 2754.4546
$ws.Range("I116").Value = 2805.1667
$ws.Range("J116").Value = 2526.25
$ws.Range("K116").Value = 2805.1667
$ws.Range("L116").Value = 2526.25
$ws.Range("M116").Value = 636.8332999999998
$ws.Range("N116").Value = -9410.25
$ws.Range("H118").Value = 492.5
$ws.Range("I118").Value = 492.5
$ws.Range("K118").Value = 1477.5
$ws.Range("M118").Value = 179.5
$ws.Range("H123").Value = 47720.555
$ws.Range("J123").Value = 47720.555
$ws.Range("L123").Value = 47720.555
$ws.Range("N123").Value = -57520.555
$ws.Range("H138").Value = 2291.9822
$ws.Range("I138").Value = 1375.8649
$ws.Range("J138").Value = 4076
$ws.Range("K138").Value = 4127.5947
$ws.Range("L138").Value = 12228
$ws.Range("M138").Value = 1012.4053
$ws.Range("N138").Value = -22508
$ws.Range("H139").Value = 11881.177
$ws.Range("J139").Value = 11881.177
$ws.Range("L139").Value = 11881.177
$ws.Range("N139").Value = -22161.177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4621.107
$ws.Range("I74").Value = 5750.7617
$ws.Range("K74").Value = 5750.7617
$ws.Range("M74").Value = -4876.7617
$ws.Range("H77").Value = 4621.107
$ws.Range("I77").Value = 5750.7617
$ws.Range("K77").Value = 28753.8085
$ws.Range("M77").Value = -24385.8085
$ws.Range("H110").Value = 2231.7646
$ws.Range("I110").Value = 2198.6667
$ws.Range("J110").Value = 2480
$ws.Range("K110").Value = 2198.6667
$ws.Range("L110").Value = 2480
$ws.Range("M110").Value = -153.6667000000002
$ws.Range("N110").Value = -6570
$ws.Range("H118").Value = 29999
$ws.Range("J118").Value = 29999
$ws.Range("L118").Value = 29999
$ws.Range("N118").Value = -33313

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 727.6799999999999
$ws.Range("I94").Value = 638.7826
$ws.Range("K94").Value = 638.7826
$ws.Range("M94").Value = -187.7826
$ws.Range("H107").Value = 2053.75
$ws.Range("J107").Value = 1993.25
$ws.Range("L107").Value = 1993.25
$ws.Range("N107").Value = -5833.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2639.1
$ws.Range("I16").Value = 2029.8334
$ws.Range("J16").Value = 3553
$ws.Range("K16").Value = 2029.8334
$ws.Range("L16").Value = 3553
$ws.Range("M16").Value = -1742.8334
$ws.Range("N16").Value = -4127
$ws.Range("H31").Value = 3652.1914
$ws.Range("I31").Value = 3023.7144
$ws.Range("K31").Value = 3023.7144
$ws.Range("M31").Value = -2728.7144
$ws.Range("H34").Value = 3652.1914
$ws.Range("I34").Value = 3023.7144
$ws.Range("K34").Value = 3023.7144
$ws.Range("M34").Value = -2821.7144
$ws.Range("H107").Value = 558.30554
$ws.Range("I107").Value = 482.95456
$ws.Range("J107").Value = 676.7143
$ws.Range("K107").Value = 482.95456
$ws.Range("L107").Value = 676.7143
$ws.Range("M107").Value = 1437.04544
$ws.Range("N107").Value = -4516.7143
$ws.Range("H113").Value = 2639.1
$ws.Range("I113").Value = 2029.8334
$ws.Range("J113").Value = 3553
$ws.Range("K113").Value = 2029.8334
$ws.Range("L113").Value = 3553
$ws.Range("M113").Value = 140.1666
$ws.Range("N113").Value = -7893

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 28444606
$ws.Range("I11").Value = 14762853
$ws.Range("K11").Value = 14762853
$ws.Range("M11").Value = -14762714
$ws.Range("H107").Value = 280.88235
$ws.Range("I107").Value = 241.5
$ws.Range("J107").Value = 375.4
$ws.Range("K107").Value = 241.5
$ws.Range("L107").Value = 375.4
$ws.Range("M107").Value = 1678.5
$ws.Range("N107").Value = -4215.4
$ws.Range("H113").Value = 7171.2383
$ws.Range("J113").Value = 14073.777
$ws.Range("L113").Value = 14073.777
$ws.Range("N113").Value = -18413.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6016.4287
$ws.Range("I61").Value = 6993.2354
$ws.Range("J61").Value = 1865
$ws.Range("K61").Value = 6993.2354
$ws.Range("L61").Value = 1865
$ws.Range("M61").Value = -6791.2354
$ws.Range("N61").Value = -2269
$ws.Range("H93").Value = 15473
$ws.Range("I93").Value = 25575.75
$ws.Range("J93").Value = 2002.6666
$ws.Range("K93").Value = 25575.75
$ws.Range("L93").Value = 2002.6666
$ws.Range("M93").Value = -24327.75
$ws.Range("N93").Value = -4498.6666
$ws.Range("H113").Value = 6016.4287
$ws.Range("I113").Value = 6993.2354
$ws.Range("J113").Value = 1865
$ws.Range("K113").Value = 6993.2354
$ws.Range("L113").Value = 1865
$ws.Range("M113").Value = -4823.2354
$ws.Range("N113").Value = -6205

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 642.4231
$ws.Range("J107").Value = 628.4167
$ws.Range("L107").Value = 1885.2501
$ws.Range("N107").Value = -5725.2501
$ws.Range("H113").Value = 514.3570999999999
$ws.Range("I113").Value = 365.1
$ws.Range("J113").Value = 887.5
$ws.Range("K113").Value = 1095.3
$ws.Range("L113").Value = 2662.5
$ws.Range("M113").Value = 1074.7
$ws.Range("N113").Value = -7002.5
